$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text value (or $null if unchanged),
# new Volume(1h) (E) text value.
$updates = @(
    @{ Row = 2;  D = "63.554.91";  E = "  +1.45%  " },
    @{ Row = 3;  D = "2.650.25";   E = "  +2.93%  " },
    @{ Row = 4;  D = $null;        E = "  -0.04%  " },
    @{ Row = 5;  D = "591.33";     E = "  +1.82%  " },
    @{ Row = 6;  D = "144.11";     E = "  -0.50%  " },
    @{ Row = 7;  D = $null;        E = "  -0.03%  " },
    @{ Row = 8;  D = $null;        E = "  -1.00%  " },
    @{ Row = 9;  D = "2.648.08";   E = "  +2.88%  " },
    @{ Row = 10; D = $null;        E = "  +0.60%  " },
    @{ Row = 11; D = $null;        E = "  +0.75%  " },
    @{ Row = 12; D = $null;        E = "  +0.80%  " },
    @{ Row = 13; D = $null;        E = "  +0.81%  " },
    @{ Row = 14; D = $null;        E = "  +1.88%  " },
    @{ Row = 15; D = "3.122.93";   E = "  +2.89%  " },
    @{ Row = 16; D = "63.470.91";  E = "  +1.46%  " },
    @{ Row = 17; D = "0.0000146";  E = "  +0.95%  " },
    @{ Row = 18; D = "2.651.52";   E = "  +3.44%  " },
    @{ Row = 19; D = "11.44";      E = "  +2.21%  " },
    @{ Row = 20; D = "340.53";     E = "  +0.87%  " },
    @{ Row = 21; D = "4.37";       E = "  +0.57%  " },
    @{ Row = 22; D = "6.74";       E = "  +1.51%  " },
    @{ Row = 23; D = $null;        E = "  +0.05%  " },
    @{ Row = 24; D = "67.19";      E = "  +0.32%  " },
    @{ Row = 25; D = $null;        E = "  +6.30%  " },
    @{ Row = 26; D = "1.54";       E = "  +4.82%  " },
    @{ Row = 27; D = $null;        E = "  +0.85%  " },
    @{ Row = 28; D = "547.63";     E = "  +18.91%  " },
    @{ Row = 29; D = "1.00";       E = "  +0.10%  " },
    @{ Row = 30; D = "8.42";       E = "  +2.53%  " },
    @{ Row = 31; D = "7.77";       E = "  -1.07%  " },
    @{ Row = 32; D = $null;        E = "  +13.95%  " },
    @{ Row = 33; D = $null;        E = "  +2.92%  " },
    @{ Row = 34; D = "0.0₃0808";   E = "  -0.29%  " },
    @{ Row = 35; D = "175.25";     E = "  -1.08%  " },
    @{ Row = 36; D = "4.90";       E = "  +9.57%  " },
    @{ Row = 37; D = $null;        E = "  -0.02%  " },
    @{ Row = 38; D = $null;        E = "  +0.54%  " },
    @{ Row = 39; D = "19.08";      E = "  +1.02%  " },
    @{ Row = 40; D = "1.81";       E = "  +7.86%  " },
    @{ Row = 41; D = $null;        E = "  -0.03%  " },
    @{ Row = 42; D = "170.20";     E = "  +7.96%  " },
    @{ Row = 43; D = "40.28";      E = "  +2.18%  " },
    @{ Row = 44; D = $null;        E = "  +0.85%  " },
    @{ Row = 45; D = "22.39";      E = "  +6.49%  " },
    @{ Row = 46; D = "0.632";      E = "  +0.84%  " },
    @{ Row = 47; D = "0.0557";     E = "  +4.30%  " },
    @{ Row = 48; D = "0.0961";     E = "  -0.22%  " },
    @{ Row = 49; D = $null;        E = "  +2.45%  " },
    @{ Row = 50; D = "18.82";      E = "  +3.97%  " },
    @{ Row = 51; D = $null;        E = "  +0.63%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force the Price cell to remain plain text (many values, e.g.
        # "591.33", would otherwise be auto-coerced to a number by Excel),
        # then drop back to the default "Normal" style so no stray
        # number-format style sticks to the cell.
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
